$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (safe from numeric auto-conversion)
$ws.Range("D2").Value = "69.638.65"
$ws.Range("E2").Value = "  -3.73%  "
$ws.Range("D3").Value = "2.511.79"
$ws.Range("E3").Value = "  -4.80%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "2.509.41"
$ws.Range("E10").Value = "  -7.57%  "
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("E14").Value = "  -4.87%  "
$ws.Range("D15").Value = "69.467.94"
$ws.Range("E15").Value = "  -3.86%  "
$ws.Range("E16").Value = "  -6.42%  "
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "2.516.53"
$ws.Range("E18").Value = "  -4.81%  "
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("E20").Value = "  -6.47%  "
$ws.Range("E21").Value = "  -7.07%  "
$ws.Range("E22").Value = "  -4.03%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("E23").Value = "  -5.58%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("E26").Value = "  -6.04%  "
$ws.Range("E27").Value = "  -5.80%  "
$ws.Range("D28").Value = "2.639.69"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").Value = "0.0₃0904"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -6.20%  "
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("E43").Value = "  -6.93%  "
$ws.Range("E44").Value = "  -14.53%  "
$ws.Range("E45").Value = "  -10.51%  "
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("E47").Value = "  -5.70%  "
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("E49").Value = "  -4.37%  "
$ws.Range("E50").Value = "  -4.58%  "
$ws.Range("E51").Value = "  -1.61%  "

# Numeric-looking text updates in column D: use a quoted-string formula then
# paste-special as values to keep them as text without altering cell style/NumberFormat
$ws.Range("D5").Formula = "=""579.07"""
$ws.Range("D6").Formula = "=""167.46"""
$ws.Range("D8").Formula = "=""0.521"""
$ws.Range("D10").Formula = "=""0.159"""
$ws.Range("D12").Formula = "=""0.346"""
$ws.Range("D17").Formula = "=""24.97"""
$ws.Range("D19").Formula = "=""7.81"""
$ws.Range("D20").Formula = "=""11.31"""
$ws.Range("D21").Formula = "=""348.39"""
$ws.Range("D22").Formula = "=""3.95"""
$ws.Range("D23").Formula = "=""1.94"""
$ws.Range("D24").Formula = "=""1.00"""
$ws.Range("D25").Formula = "=""68.76"""
$ws.Range("D27").Formula = "=""8.95"""
$ws.Range("D29").Formula = "=""0.997"""
$ws.Range("D32").Formula = "=""1.27"""
$ws.Range("D33").Formula = "=""465.26"""
$ws.Range("D37").Formula = "=""154.28"""
$ws.Range("D38").Formula = "=""18.97"""
$ws.Range("D39").Formula = "=""18.38"""
$ws.Range("D41").Formula = "=""4.76"""
$ws.Range("D45").Formula = "=""2.30"""
$ws.Range("D46").Formula = "=""38.08"""
$ws.Range("D47").Formula = "=""142.90"""
$ws.Range("D48").Formula = "=""0.530"""
$ws.Range("D50").Formula = "=""1.60"""
$ws.Range("D51").Formula = "=""0.0733"""

$ws.Range("D2:D51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$excel.CutCopyMode = $false

